# readExcelForScan: locate the "d" cell in row 9 (B9) and replace it with "scan",
# formatted with an orange Consolas font, vertically centered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("B9")
$target.Value = "scan"

$target.Font.Name = "Consolas"
$target.Font.Color = 7901646   # 0xCE9178 as BGR long (RGB 206,145,120)
$target.VerticalAlignment = -4108   # xlCenter

$target.Select()
